# Apply the commit "Added 'Date and Time' and 'Cycle_count' parameters"
# to the Analysis Results worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Structural changes: insert/delete rows so that every remaining
#    (unchanged-value) row lands on its correct final row number.
# ---------------------------------------------------------------------

# Insert a new row 1 for "Date and Time" - this also shifts the
# [hh]:mm:ss number format that lived on B1 down to B2, which is what
# we want since "Total time taken for the ride" moves to row 2.
$ws.Rows.Item(1).Insert()

# The old "Maximum BMS Temperature in C" row (now at row 32) is removed
# entirely in the new version - no equivalent row remains.
$ws.Rows.Item(32).Delete()

# Insert a new row for "Cycle Count of battery" right before the row
# that will hold "Idling time percentage" (currently row 35).
$ws.Rows.Item(35).Insert()

# ---------------------------------------------------------------------
# 2) Write the full final A:B content for rows 1-45.
# ---------------------------------------------------------------------

function Set-Row([int]$r, [string]$label, $value, [bool]$isText) {
    $ws.Cells.Item($r, 1).Value = $label
    if ($isText) {
        $ws.Cells.Item($r, 2).Value = $value
    } else {
        $ws.Cells.Item($r, 2).Value = [double]$value
    }
}

Set-Row 1  "Date and Time" "2024-03-12 17:38:19.949000 to 2024-03-12 19:29:39.538000" $true
Set-Row 2  "Total time taken for the ride" "0.07714497685185186" $false
Set-Row 3  "Actual Ampere-hours (Ah)" "31.12770305555556" $false
Set-Row 4  "Actual Watt-hours (Wh)" "1615.390856283611" $false
Set-Row 5  "Starting SoC (Ah)" "5.944" $false
Set-Row 6  "Ending SoC (Ah)" "6.418" $false
Set-Row 7  "Starting SoC (%)" "15" $false
Set-Row 8  "Ending SoC (%)" "100" $false
Set-Row 9  "Total distance covered (km)" "57.7466072884102" $false
Set-Row 10 "Total energy consumption(WH/KM)" "27.97377945020541" $false
Set-Row 11 "Total SOC consumed(%)" "85" $false
$ws.Cells.Item(12, 1).Value = "Mode"
$ws.Cells.Item(12, 2).Value = "Eco mode`n85.95%`nSports mode`n13.67%"
Set-Row 13 "Peak Power(kW)" "4429.990945" $false
Set-Row 14 "Average Power(kW)" "-877.0030285961393" $false
Set-Row 15 "Total Energy Regenerated(kWh)" "68.13299414972222" $false
Set-Row 16 "Regenerative Effectiveness(%)" "4.047046564394381" $false
Set-Row 17 "Highest Cell Voltage(V)" "3.483" $false
Set-Row 18 "Lowest Cell Voltage(V)" "3.144" $false
Set-Row 19 "Difference in Cell Voltage(V)" "0.339" $false
Set-Row 20 "Minimum Temperature(C)" "36" $false
Set-Row 21 "Maximum Temperature(C)" "48" $false
Set-Row 22 "Difference in Temperature(C)" "12" $false
Set-Row 23 "Maximum Fet Temperature-BMS(C)" "58" $false
Set-Row 24 "Maximum Afe Temperature-BMS(C)" "63" $false
Set-Row 25 "Maximum PCB Temperature-BMS(C)" "59" $false
Set-Row 26 "Maximum MCU Temperature(C)" "43" $false
Set-Row 27 "Maximum Motor Temperature(C)" "0" $false
Set-Row 28 "Abnormal Motor Temperature Detected(C)" "0" $false
Set-Row 29 "highest cell temp(C)" "48" $false
Set-Row 30 "lowest cell temp(C)" "36" $false
Set-Row 31 "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)" "12" $false
Set-Row 32 "Battery Voltage(V)" "56" $false
Set-Row 33 "Total energy charged(kWh)" "1.743151371111111" $false
Set-Row 34 "Electricity consumption units(kW)" "7.264946949700389e-08" $false
Set-Row 35 "Cycle Count of battery" "138" $false
Set-Row 36 "Idling time percentage" "7.10488322717622" $false
Set-Row 37 "Time spent in 0-10 km/h" "5.470912951167728" $false
Set-Row 38 "Time spent in 10-20 km/h" "8.798301486199575" $false
Set-Row 39 "Time spent in 20-30 km/h" "14.21656050955414" $false
Set-Row 40 "Time spent in 30-40 km/h" "41.97537154989384" $false
Set-Row 41 "Time spent in 40-50 km/h" "21.1176220806794" $false
Set-Row 42 "Time spent in 50-60 km/h" "0" $false
Set-Row 43 "Time spent in 60-70 km/h" "0" $false
Set-Row 44 "Time spent in 70-80 km/h" "0" $false
Set-Row 45 "Time spent in 80-90 km/h" "0" $false

# ---------------------------------------------------------------------
# 3) Make sure the "Total time taken for the ride" duration cell keeps
#    its [hh]:mm:ss display format (it should already, from the Insert
#    at step 1, but set it explicitly to be safe).
# ---------------------------------------------------------------------
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

Write-Host "Done."
